{"js": "// Remove the trailing \"Ver no Jupiter...\" / copyright footer block (and the\n// blank paragraph preceding it) that followed the last Bibliografia entry.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraph (\"... Thomson Pioneira (2008).\") that must stay,\n// and the two footer paragraphs that must go away, by exact text match so the\n// edit is robust to any positional drift.\nconst anchorText =\n  \"(2009).JEWETT Jr, John W.; SERWAY, Raymond A. Princ\u00edpios de F\u00edsica. Vol. 1, Thomson Pioneira (2008).\";\nconst jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst copyrightText =\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\nlet anchorIndex = -1;\nlet jupiterIndex = -1;\nlet copyrightIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (t === anchorText) anchorIndex = i;\n  else if (t === jupiterText) jupiterIndex = i;\n  else if (t === copyrightText) copyrightIndex = i;\n}\n\nif (jupiterIndex !== -1 && copyrightIndex !== -1) {\n  // The empty paragraph that sits right before the \"Ver no Jupiter...\"\n  // paragraph is also removed by the edit (it only makes sense as a spacer\n  // before the footer block that is going away).\n  let blankIndex = jupiterIndex - 1;\n  if (anchorIndex !== -1 && blankIndex === anchorIndex) blankIndex = -1;\n  if (blankIndex !== -1 && items[blankIndex].text !== \"\") blankIndex = -1;\n\n  // Delete highest index first so the other queued proxies stay valid.\n  const toDelete = [copyrightIndex, jupiterIndex];\n  if (blankIndex !== -1) toDelete.push(blankIndex);\n  toDelete.sort((a, b) => b - a);\n  for (const idx of toDelete) items[idx].delete();\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / copyright footer block (and the\n# blank paragraph preceding it) that followed the last Bibliografia entry.\n$d = $word.ActiveDocument\n\n$anchorText = \"(2009).JEWETT Jr, John W.; SERWAY, Raymond A. Princ\u00edpios de F\u00edsica. Vol. 1, Thomson Pioneira (2008).\"\n$jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$anchorIndex = -1\n$jupiterIndex = -1\n$copyrightIndex = -1\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)\n  if ($t -eq $anchorText) { $anchorIndex = $i }\n  elseif ($t -eq $jupiterText) { $jupiterIndex = $i }\n  elseif ($t -eq $copyrightText) { $copyrightIndex = $i }\n}\n\nif ($jupiterIndex -ne -1 -and $copyrightIndex -ne -1) {\n  # The empty paragraph right before the \"Ver no Jupiter...\" paragraph is\n  # also removed (it only served as a spacer before the footer block).\n  $blankIndex = $jupiterIndex - 1\n  if ($anchorIndex -ne -1 -and $blankIndex -eq $anchorIndex) { $blankIndex = -1 }\n  if ($blankIndex -ne -1 -and $d.Paragraphs.Item($blankIndex).Range.Text.TrimEnd([char]13) -ne \"\") { $blankIndex = -1 }\n\n  # Delete from the highest index down so earlier indices stay valid.\n  $indices = @($copyrightIndex, $jupiterIndex)\n  if ($blankIndex -ne -1) { $indices += $blankIndex }\n  $indices = $indices | Sort-Object -Descending\n\n  foreach ($idx in $indices) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n  }\n}\n"}
